# Adds upstream emissions/removals data for alloy variants of BOF use.
#
# Semantics (per commit message "added upstream for alloy"):
#   - "BOF use" (generic) is renamed to "BOF use + low alloy" on the row that used
#     to hold the generic BOF figures, and a NEW generic "BOF use" row is introduced
#     using the values that used to live under "BOF use - no alloy".
#   - Three new substance rows are appended: "charcoal kiln use" (label only, no data
#     yet), "no alloy", "chromium alloy" and "low alloy" factory-use coefficients.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("emissions")
$ws2 = $wb.Worksheets.Item("removals")

# ---------------------------------------------------------------------------
# Sheet "emissions"
# ---------------------------------------------------------------------------

# Existing "BOF use" row becomes the low-alloy variant.
$ws1.Range("A50").Value = "BOF use + low alloy"

# Existing "BOF use - no alloy" row becomes the plain "BOF use" row.
$ws1.Range("A57").Value = "BOF use"

# New row 62: label only.
$ws1.Range("A62").Value = "charcoal kiln use"

# New row 63: "no alloy"
$ws1.Range("A63").Value = "no alloy"
$ws1.Range("B63").Formula = "=SUM(C63:E63)"
$ws1.Range("C63").Value = [double]"5.1999999999999998E-2"
$ws1.Range("D63").Value = [double]"2.8700000000000002E-3"
$ws1.Range("E63").Value = [double]"3.8800000000000001E-5"
$ws1.Range("E63").NumberFormat = "0.00E+00"
$ws1.Range("G63").Formula = "=SUM(H63:J63)"
$ws1.Range("H63").Formula = "=0.00572"
$ws1.Range("I63").Formula = "=0.0000465"
$ws1.Range("J63").Value = [double]"1.2200000000000001E-7"
$ws1.Range("J63").NumberFormat = "0.00E+00"

# New row 64: "chromium alloy"
$ws1.Range("A64").Value = "chromium alloy"
$ws1.Range("B64").Formula = "=SUM(C64:E64)"
$ws1.Range("C64").Value = [double]"3.33"
$ws1.Range("D64").Value = [double]"0.182"
$ws1.Range("E64").Value = [double]"2.64E-3"
$ws1.Range("G64").Formula = "=SUM(H64:J64)"
$ws1.Range("H64").Value = [double]"0.36699999999999999"
$ws1.Range("I64").Value = [double]"2.9299999999999999E-3"
$ws1.Range("J64").Value = [double]"8.1499999999999999E-6"
$ws1.Range("J64").NumberFormat = "0.00E+00"

# New row 65: "low alloy"
$ws1.Range("A65").Value = "low alloy"
$ws1.Range("B65").Formula = "=SUM(C65:E65)"
$ws1.Range("C65").Value = [double]"0.49099999999999999"
$ws1.Range("D65").Value = [double]"2.6200000000000001E-2"
$ws1.Range("E65").Value = [double]"4.4999999999999999E-4"
$ws1.Range("G65").Formula = "=SUM(H65:J65)"
$ws1.Range("H65").Value = [double]"5.3900000000000003E-2"
$ws1.Range("I65").Value = [double]"4.9200000000000003E-4"
$ws1.Range("J65").Value = [double]"1.8500000000000001E-6"
$ws1.Range("J65").NumberFormat = "0.00E+00"

# ---------------------------------------------------------------------------
# Sheet "removals"
# ---------------------------------------------------------------------------

# Existing "BOF use" row becomes the low-alloy variant.
$ws2.Range("A49").Value = "BOF use - low alloy"

# Existing "BOF use - no alloy" row becomes the plain "BOF use" row.
$ws2.Range("A56").Value = "BOF use"

# New rows 61-63.
$ws2.Range("A61").Value = "no alloy"
$ws2.Range("B61").Value = [double]"2.5500000000000002E-3"

$ws2.Range("A62").Value = "chromium alloy"
$ws2.Range("B62").Value = [double]"0.16"

$ws2.Range("A63").Value = "low alloy"
$ws2.Range("B63").Value = [double]"2.3699999999999999E-2"
